$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 2014
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 27
$ws.Range("D29").Value = 0.34375
$ws.Range("E29").Value = 0.54166666666666663

$ws.Range("F28").AutoFill($ws.Range("F28:F29"))
$ws.Range("G28").AutoFill($ws.Range("G28:G29"))

$ws.Range("F29").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
